# "Forgot to save excel" - persist the pending workbook changes:
#  - add a new "IS" sheet (GAM / DICKEBL / GAMBL lookup table) after "DN"
#  - make it the active sheet/tab (zoom 145%, selection on D3)
#  - the previously active sheet ("DN") stops being the selected tab

$wb = $excel.ActiveWorkbook

# Add the new worksheet as the LAST tab (after the current last sheet, "DN")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "IS"

# Pre-format the 4 used columns as Text so the numeric-looking values
# (130 / 1 / 2700) are stored as text, same as the rest of the workbook.
$ws.Range("A1:D2").NumberFormat = "@"

# Fill in column-by-column (matches shared-string authoring order)
$ws.Range("B1").Value = "GAM"
$ws.Range("B2").Value = "130"
$ws.Range("C1").Value = "DICKEBL"
$ws.Range("C2").Value = "1"
$ws.Range("D1").Value = "GAMBL"
$ws.Range("D2").Value = "2700"
$ws.Range("A1").Value = "NAME"
$ws.Range("A2").Value = "ISOLERING"

# The other sheets in this workbook carry a (blank/grey) sheet background
# picture; mirror that on the new sheet too (best effort - harmless if the
# host can't resolve/persist an external image file).
try {
    $ws.SetBackgroundPicture("image4.png")
} catch {
}

# Match the view state recorded in the saved file: zoomed to 145%, the
# active/selected cell on the new sheet is D3.
[void]$ws.Range("D3").Select()
$excel.ActiveWindow.Zoom = 145

Write-Host "IS sheet added and populated"
